# Update column G ("K", formerly "Strike#") with regenerated strikeout
# values for Drew Pomeranz's 2021 game log. These values are recomputed
# from the source data (s_vals) and written back as literal numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 2
    6  = 2
    7  = 1
    8  = 2
    9  = 0
    10 = 2
    11 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 2
    18 = 0
    19 = 1
    20 = 3
    21 = 1
    22 = 3
    23 = 2
    24 = 1
    25 = 3
    26 = 1
    27 = 0
    28 = 3
    29 = 3
    30 = 2
    31 = 2
    32 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
